$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new "Person" trivia rows below the existing two.
# Cells are written in the same order the original author apparently used
# (answers for the first two new rows, then their questions, then the
# third row's answer and question) so the workbook's shared-strings table
# ends up built in the same order as the target file.
$ws.Range("B4").Value = "Elon Musk"
$ws.Range("B5").Value = "Lionel Messi"
$ws.Range("A5").Value = "Who is the best football player?"
$ws.Range("A4").Value = "Who is CEO of SpaceX?"
$ws.Range("B6").Value = "Albert Einstein"
$ws.Range("A6").Value = "Who is most famous for inventing the theory of relativity?"

$ws.Range("C4").Value = "Person"
$ws.Range("C5").Value = "Person"
$ws.Range("C6").Value = "Person"

# Column A needs to widen to fit the new, longer questions (closest
# achievable width to the recorded best-fit value).
$ws.Columns.Item(1).ColumnWidth = 52.5

# Leave the selection where the author apparently left it when saving.
$ws.Range("P15").Select() | Out-Null
